# Portland Home Miscellaneous - add a "Match ID" column at the front of
# the sheet (multiple cleaning changes, added sql schema).
#
# Net effect vs. the previous layout: every existing column shifts one
# slot to the right (A->B, B->C, ... W->X) and the new column A is
# populated with the match id (10) for every player / totals row, with
# a header label of "Match ID" in row 2 (the visible header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the whole sheet one column to the right to make room for the
# new "Match ID" column.
$ws.Columns.Item(1).Insert()

# Rows 1, 3 and 20 are hidden; temporarily reveal them so writing into
# them does not perturb their row height, then restore hidden state.
$ws.Rows.Item(1).Hidden = $false
$ws.Rows.Item(3).Hidden = $false
$ws.Rows.Item(20).Hidden = $false

# Header (row 2) label for the new column.
$ws.Range("A2").Value = "Match ID"

# Data + totals rows: the match id is constant (10) for every row of
# this sheet.
for ($r = 4; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = 10
}

# Match the look of the other header cells: bold text, no border, no
# special alignment (row 20's cell intentionally keeps the default
# style, matching the rest of that hidden totals row).
$ws.Range("A2:A19").Font.Bold = $true

# Restore hidden rows.
$ws.Rows.Item(1).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(20).Hidden = $true

# Selection moves onto the freshly filled column.
[void]$ws.Range("A2:A19").Select()
